$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Info")

# Update values: Main_IP+mask (B12) from "dhcp" to a real static IP/mask
$ws.Range("B12").Value = "1.1.1.1/24"

# Backup_GW (B21) was empty, now has a gateway value
$ws.Range("B21").Value = "1.1.1.12"

# Backup_IP+mask (B20) from "DHCP" to a real static IP/mask
$ws.Range("B20").Value = "2.2.2.2/24"

# Update the view's active cell / scroll position
$ws.Application.Goto($ws.Range("A7"), $true)
$ws.Range("G16").Select()
